$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (BGR-packed integers, as consumed by Range.Interior.Color)
$yellow = 65535      # RGB(255,255,0)
$red    = 255        # RGB(255,0,0)
$blue   = 13998939   # RGB(91,155,213) -- same accent color already used elsewhere (theme accent5)

# --- Header row: replace "RFR Final" / "KNR Final" labels with a 1st/2nd/3rd ranking ---
$ws.Range("F1").ClearContents()
$ws.Range("G1").Value = "1st"
$ws.Range("G1").Interior.Color = $blue

$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "2nd"
$ws.Range("G2").Interior.Color = $yellow

$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "3rd"
$ws.Range("G3").Interior.Color = $red

# --- Highlight the best-performing model per remaining row (yellow = 2nd, red = 3rd place) ---
$ws.Range("C4").Interior.Color = $yellow
$ws.Range("C5").Interior.Color = $yellow
$ws.Range("D6").Interior.Color = $yellow
$ws.Range("C7").Interior.Color = $yellow

$ws.Range("C6").Interior.Color = $red
$ws.Range("C8").Interior.Color = $red
$ws.Range("E8").Interior.Color = $red

# --- Drop the retired RFR/KNR prediction columns (F is cleared, G is removed outright) ---
$ws.Range("F4").ClearContents()
$ws.Range("G4").Clear()

$ws.Range("F5").ClearContents()
$ws.Range("G5").Clear()

$ws.Range("F6").ClearContents()
$ws.Range("G6").Clear()

$ws.Range("F7").ClearContents()
$ws.Range("G7").Clear()

$ws.Range("F8").ClearContents()
$ws.Range("G8").Clear()

# Column F width (was auto; now fixed, matching the new narrower "empty" column)
$ws.Columns.Item(6).ColumnWidth = 10

# Move the active selection
$ws.Range("K12").Select()
